$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date in column C for rows 2-8
# from 2023-10-09 (45208) to 2023-10-13 (45212)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 45212
}
